# Applies updated pl_mw.xlsx res_line values for Case_4_183 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.25442998338707
$ws.Range("C2").Value = 0.305907168253583
$ws.Range("D2").Value = 0.07910933540964038
$ws.Range("E2").Value = 0.4189847875735495
$ws.Range("G2").Value = 1.320207317101534
$ws.Range("H2").Value = 1.060871089138544
# Row 3
$ws.Range("B3").Value = 1.11571172364205
$ws.Range("C3").Value = 0.2660962492269334
$ws.Range("D3").Value = 0.07163942886302266
$ws.Range("E3").Value = 0.3650424495941564
$ws.Range("G3").Value = 1.234646230428154
$ws.Range("H3").Value = 1.027484546662464
# Row 4
$ws.Range("B4").Value = 1.031038508535062
$ws.Range("C4").Value = 0.2416618522267129
$ws.Range("D4").Value = 0.06709830182005305
$ws.Range("E4").Value = 0.3320775187322766
$ws.Range("G4").Value = 1.183124712801799
$ws.Range("H4").Value = 1.007731740605379
# Row 5
$ws.Range("B5").Value = 0.996656222508193
$ws.Range("C5").Value = 0.231706091892022
$ws.Range("D5").Value = 0.06525890993091821
$ws.Range("E5").Value = 0.3186797070689806
$ws.Range("G5").Value = 1.16237764516157
$ws.Range("H5").Value = 0.9998671709435882
# Row 6
$ws.Range("B6").Value = 0.9909543949738122
$ws.Range("C6").Value = 0.2300530054021692
$ws.Range("D6").Value = 0.06495414624363605
$ws.Range("E6").Value = 0.316457062536756
$ws.Range("G6").Value = 1.158947423993141
$ws.Range("H6").Value = 0.9985723460179656
# Row 7
$ws.Range("B7").Value = 1.030574324050008
$ws.Range("C7").Value = 0.2415275807410069
$ws.Range("D7").Value = 0.06707345032654644
$ws.Range("E7").Value = 0.3318966913833918
$ws.Range("G7").Value = 1.182843913426439
$ws.Range("H7").Value = 1.007624931584331
# Row 8
$ws.Range("B8").Value = 1.206494151979427
$ws.Range("C8").Value = 0.2921777489528381
$ws.Range("D8").Value = 0.07652411837329964
$ws.Range("E8").Value = 0.4003509593017895
$ws.Range("G8").Value = 1.290491617356878
$ws.Range("H8").Value = 1.049202676127322
# Row 9
$ws.Range("B9").Value = 1.555603458338112
$ws.Range("C9").Value = 0.3916293205220995
$ws.Range("D9").Value = 0.0954305089502725
$ws.Range("E9").Value = 0.5359937010894669
$ws.Range("G9").Value = 1.509931604579435
$ws.Range("H9").Value = 1.136797315845286
# Row 10
$ws.Range("B10").Value = 1.814867178282213
$ws.Range("C10").Value = 0.4648516556011373
$ws.Range("D10").Value = 0.1095687659543216
$ws.Range("E10").Value = 0.6367497022969246
$ws.Range("G10").Value = 1.676677531922593
$ws.Range("H10").Value = 1.205045105801844
# Row 11
$ws.Range("B11").Value = 1.933473796554665
$ws.Range("C11").Value = 0.4982140119613518
$ws.Range("D11").Value = 0.1160589031314316
$ws.Range("E11").Value = 0.6828792382898001
$ws.Range("G11").Value = 1.753830512980301
$ws.Range("H11").Value = 1.236981555098112
# Row 12
$ws.Range("B12").Value = 1.978487233861586
$ws.Range("C12").Value = 0.510856471565944
$ws.Range("D12").Value = 0.1185253159803068
$ws.Range("E12").Value = 0.7003940924026182
$ws.Range("G12").Value = 1.783240870065754
$ws.Range("H12").Value = 1.24920654137992
# Row 13
$ws.Range("B13").Value = 1.968788287986229
$ws.Range("C13").Value = 0.5081332778456158
$ws.Range("D13").Value = 0.1179937360783043
$ws.Range("E13").Value = 0.6966198153732392
$ws.Range("G13").Value = 1.7768980683216
$ws.Range("H13").Value = 1.24656777644401
# Row 14
$ws.Range("B14").Value = 1.937175058583989
$ws.Range("C14").Value = 0.4992539300468479
$ws.Range("D14").Value = 0.1162616394269094
$ws.Range("E14").Value = 0.6843192377623808
$ws.Range("G14").Value = 1.75624618628126
$ws.Range("H14").Value = 1.237984660390737
# Row 15
$ws.Range("B15").Value = 1.917824142181757
$ws.Range("C15").Value = 0.4938162637143364
$ws.Range("D15").Value = 0.1152018277609415
$ws.Range("E15").Value = 0.6767909726522134
$ws.Range("G15").Value = 1.743621822354669
$ws.Range("H15").Value = 1.232744459117839
# Row 16
$ws.Range("B16").Value = 1.807129691975717
$ws.Range("C16").Value = 0.4626725196745269
$ws.Range("D16").Value = 0.1091458257185565
$ws.Range("E16").Value = 0.6337413208624696
$ws.Range("G16").Value = 1.671662175906903
$ws.Range("H16").Value = 1.202976180685823
# Row 17
$ws.Range("B17").Value = 1.739395591814628
$ws.Range("C17").Value = 0.443581306432236
$ws.Range("D17").Value = 0.1054458892147636
$ws.Range("E17").Value = 0.6074102734886395
$ws.Range("G17").Value = 1.627855265210542
$ws.Range("H17").Value = 1.184944662353132
# Row 18
$ws.Range("B18").Value = 1.700499199300168
$ws.Range("C18").Value = 0.4326054727962969
$ws.Range("D18").Value = 0.10332327745013
$ws.Range("E18").Value = 0.5922928754959997
$ws.Range("G18").Value = 1.602780550790669
$ws.Range("H18").Value = 1.174656995709938
# Row 19
$ws.Range("B19").Value = 1.687340162536202
$ws.Range("C19").Value = 0.42889005401247
$ws.Range("D19").Value = 0.1026055316425527
$ws.Range("E19").Value = 0.5871789863538197
$ws.Range("G19").Value = 1.594311391377346
$ws.Range("H19").Value = 1.171188023648398
# Row 20
$ws.Range("B20").Value = 1.746599509297141
$ws.Range("C20").Value = 0.4456130800400615
$ws.Range("D20").Value = 0.1058391829128169
$ws.Range("E20").Value = 0.6102103799578913
$ws.Range("G20").Value = 1.632505920566416
$ws.Range("H20").Value = 1.186855473446172
# Row 21
$ws.Range("B21").Value = 1.946457893353283
$ws.Range("C21").Value = 0.5018617593897829
$ws.Range("D21").Value = 0.1167701584790422
$ws.Range("E21").Value = 0.6879309172300481
$ws.Range("G21").Value = 1.762306815094178
$ws.Range("H21").Value = 1.240502138058901
# Row 22
$ws.Range("B22").Value = 2.077659840819194
$ws.Range("C22").Value = 0.5386757569559109
$ws.Range("D22").Value = 0.1239652708920573
$ws.Range("E22").Value = 0.7389994021384325
$ws.Range("G22").Value = 1.848274625830186
$ws.Range("H22").Value = 1.276330508877152
# Row 23
$ws.Range("B23").Value = 2.007580235491616
$ws.Range("C23").Value = 0.5190222411223999
$ws.Range("D23").Value = 0.1201203231432117
$ws.Range("E23").Value = 0.7117167541878615
$ws.Range("G23").Value = 1.802285639974116
$ws.Range("H23").Value = 1.257136919106074
# Row 24
$ws.Range("B24").Value = 1.743342478814952
$ws.Range("C24").Value = 0.4446945154101059
$ws.Range("D24").Value = 0.1056613607420758
$ws.Range("E24").Value = 0.6089443877483518
$ws.Range("G24").Value = 1.630403015769303
$ws.Range("H24").Value = 1.185991350267244
# Row 25
$ws.Range("B25").Value = 1.460692016564508
$ws.Range("C25").Value = 0.3647035457725565
$ws.Range("D25").Value = 0.09027378889304316
$ws.Range("E25").Value = 0.499123579319658
$ws.Range("G25").Value = 1.449629249400004
$ws.Range("H25").Value = 1.112431583821291

Write-Output "Updated 144 cells in Sheet1 (B2:H25, excluding F and I-O)"
